# Update EUR->ARS rate: 2025-10-07T15:23:17Z
# Appends a new row to the quote history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 63

# Column A holds a date-looking string ("2025-10-07"). The sheet stores it
# as plain text (inline string), not as a real date serial, so force the
# cell to be treated as text before assigning the value (mirrors typing
# into a cell that's pre-formatted as Text in Excel). ClearFormats()
# afterwards drops the temporary "Text" number format so the cell's style
# matches the rest of the sheet (no explicit style override).
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "2025-10-07"
$ws.Range("A$newRow").ClearFormats()

$ws.Range("B$newRow").Value = "15:23:17"
$ws.Range("C$newRow").Value = "1.00 EUR = 1,778.9258"
